# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (only used by the notes master)
#   ppt/theme/theme2.xml -> "Integral"     (used by the slide master / every slide)
#
# The target revision swaps the two theme colour palettes: the theme that
# backs the slides (theme2.xml) becomes the "Office" palette, while the
# notes-only theme (theme1.xml) becomes the "Integral" palette. Font scheme
# and format scheme are identical between the two themes already, so the
# only thing that actually needs to move is the 12-slot colour scheme.
#
# Re-point every slide's theme colour scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) at the "Office Theme" palette - this is the
# palette the live slide theme (theme2.xml) needs to end up with.

$p = $ppt.ActivePresentation
$slides = $p.Slides.Range()
$themeColors = $slides.ThemeColorScheme

# OLE/VBA colours are packed as 0x00BBGGRR, i.e. RGB() byte order reversed.
# Index -> (scheme slot, target hex, packed BGR integer)
#  1 -> dk1      000000 -> 0
#  2 -> lt1      FFFFFF -> 16777215
#  3 -> dk2      44546A -> 6968388
#  4 -> lt2      E7E6E6 -> 15132391
#  5 -> accent1  5B9BD5 -> 13998939
#  6 -> accent2  ED7D31 -> 3243501
#  7 -> accent3  A5A5A5 -> 10855845
#  8 -> accent4  FFC000 -> 49407
#  9 -> accent5  4472C4 -> 12874308
# 10 -> accent6  70AD47 -> 4697456
# 11 -> hlink    0563C1 -> 12673797
# 12 -> folHlink 954F72 -> 7491477
$officeThemePalette = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemePalette[$i - 1]
}
